# Applies the author's edit:
#  - Slide 16: four tables ("Table 11", "Table 13", "Table 14", "Table 15")
#    each have a header cell (row 1, col 2) whose text changes from
#    "MWh" to "Wh".
#  - Slide 17: the "Rectangle 5" textbox's second run changes its wording
#    (the shape uses spAutoFit, so its height adjusts automatically).

$p = $ppt.ActivePresentation

# --- Slide 16: MWh -> Wh in the four summary tables -----------------------
$s16 = $p.Slides.Item(16)

$tableShapeNames = @("Table 11", "Table 13", "Table 14", "Table 15")
foreach ($shapeName in $tableShapeNames) {
    $tblShape = $s16.Shapes.Item($shapeName)
    $tbl = $tblShape.Table
    $cell = $tbl.Cell(1, 2)
    $cellRange = $cell.Shape.TextFrame.TextRange
    if ($cellRange.Text -eq "MWh") {
        $cellRange.Text = "Wh"
    }
}

# --- Slide 17: update the TMY explanation sentence -------------------------
$s17 = $p.Slides.Item(17)
$rectShape = $s17.Shapes.Item("Rectangle 5")
$tr = $rectShape.TextFrame.TextRange
$run2 = $tr.Runs(2, 1)
$run2.Text = ": Typical Meteorological Year, which means it assumes the same variability of solar output for other years."
